$d = $word.ActiveDocument

# 1) Fix the grammar error "who's" -> "whose" in the sentence about the
#    sample's environment. The match is extended to include the following
#    space so that the Find/Replace fully consumes the run carrying the
#    mis-typed word (and its spell/grammar proofing marks), letting the
#    proofing marks be dropped since the flagged text no longer exists.
$d.Content.Find.Execute("people who’s ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "people whose ", 2)

# 2) The source document splits "whose" across two runs ("whos" + "e").
#    Recreate that split by carving out the single character "e" into its
#    own run. Adding then immediately removing a bookmark around that
#    character is a side-effect-free way to force Word to break the run
#    there without leaving any residual formatting behind.
$rng = $d.Content
$rng.Find.Execute("whos", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$splitPoint = $rng.End

$eRange = $d.Range($splitPoint, $splitPoint + 1)
$bm = $d.Bookmarks.Add("__split_marker__", $eRange)
$d.Bookmarks("__split_marker__").Delete()
